$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the now-unused trailing columns F:L entirely (also shrinks the
# dimension / shared-string table automatically, no leftover formatting).
$ws.Range("F1:L6").Clear()

# --- Header row (A1:C1 already correct; fix D1/E1) ---
$ws.Range("D1").Value = "SE: `$\hat\lambda_{SPF}`$(Q)"
$ws.Range("E1").Value = "SE: `$\hat\lambda_{SCE}`$(M)"

# --- Row 2 ---
$ws.Range("A2").Value = "Forecast"
$ws.Range("D2").Value = 0.04
$ws.Range("E2").Value = 1

# --- Row 3 ---
$ws.Range("A3").Value = "FE"
$ws.Range("D3").Value = 0.05
$ws.Range("E3").Value = 0.19

# --- Row 4 ---
$ws.Range("A4").Value = "FE"
$ws.Range("B4").Value = "Disg"
$ws.Range("D4").Value = 0.17
$ws.Range("E4").Value = 0.19

# --- Row 5 ---
$ws.Range("A5").Value = "FE"
$ws.Range("B5").Value = "Var"
$ws.Range("C5").Clear()
$ws.Range("D5").Value = 0.16
$ws.Range("E5").Value = 0.18

# --- Row 6 ---
$ws.Range("A6").Value = "FE"
$ws.Range("B6").Value = "Disg"
$ws.Range("C6").Value = "Var"
$ws.Range("D6").Value = 0.53
$ws.Range("E6").Value = 0.18
